$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

# New column AT holds the 29-jul prices, immediately after the existing
# 28-jul column (AS). Copy AS1's formatting (bold/centered/bordered header
# style) onto AT1, then set the header + data values for the new column.
$ws.Range("AS1").Copy($ws.Range("AT1"))

$ws.Range("AT1").Value = "29-jul"
$ws.Range("AT2").Value = 50.68
$ws.Range("AT3").Value = 37.97
$ws.Range("AT4").Value = 34.78
$ws.Range("AT5").Value = 32
$ws.Range("AT6").Value = 22.98
$ws.Range("AT7").Value = 34.78
$ws.Range("AT8").Value = 45
$ws.Range("AT9").Value = 50
$ws.Range("AT10").Value = 49.45
$ws.Range("AT11").Value = 35.61
$ws.Range("AT12").Value = 11.73
$ws.Range("AT13").Value = 25.25
$ws.Range("AT14").Value = 25.92
$ws.Range("AT15").Value = 5.34
$ws.Range("AT16").Value = 3.78
$ws.Range("AT17").Value = 8.51
$ws.Range("AT18").Value = 22.87
$ws.Range("AT19").Value = 51.53
$ws.Range("AT20").Value = 48.97
$ws.Range("AT21").Value = 72.14
$ws.Range("AT22").Value = 57.58
$ws.Range("AT23").Value = 98
$ws.Range("AT24").Value = 103.27
$ws.Range("AT25").Value = 84.13
